# "Generate Report for Handoff"
#
# The localization status report moves from "In Translation" to
# "Ready for handoff". Every place that shows that status string gets
# updated, along with the timestamps that were refreshed when the report
# was regenerated. The Status columns also get a bit wider so the new
# (longer) label isn't truncated.

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$overviewTime   = "2016-09-01 03:09:33"
$handoffTime    = "2016-09-01 03:09:29"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = $overviewTime

# Widen the zh-cn / de-de status columns to fit the new label.
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = $handoffTime
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = $overviewTime
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
